# Append a new row of log data at the bottom of the sheet
# (row 8: "Append row at 2025-05-01T11:18:34.545Z")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8
$ws.Range("A$row").Value = "2025-05-01T11:18:34.545Z"
$ws.Range("B$row").Value = "IDRF"
$ws.Range("C$row").Value = "C3"
$ws.Range("D$row").Value = "الرحلة 1"
$ws.Range("E$row").Value = "الصمود"
$ws.Range("F$row").Value = "يامن "
# G/H look like numbers ("123123" / "123"); lead with an apostrophe so they
# stay text cells, matching every other quantity/time column in this sheet.
$ws.Range("G$row").Value = "'123123"
$ws.Range("H$row").Value = "'123"

# Clear the quote-prefix style COM applies for the apostrophe-led values above
# so the new row's cells carry the same default styling as the rest of the sheet.
$ws.Range("A$($row):H$($row)").Style = "Normal"
